$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F8").Value = 4
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 1
$ws.Range("F22").Value = -1
$ws.Range("F27").Value = -4
$ws.Range("F29").Value = -4
$ws.Range("F31").Value = -6
$ws.Range("F34").Value = -1
$ws.Range("F36").Value = -2
$ws.Range("F37").Value = -2
$ws.Range("F42").Value = 5
$ws.Range("F44").Value = -1
